# 5-60.xlsx edit: add "Statement of Reconciliation" table (5-58) next to the
# existing Statement of Cash Flow, then drop the old empty
# "Reconciliation Schedule" sheet and rename the remaining sheet.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# The old workbook has two sheets: "Statement of Cash Flow" (has all the data)
# and "Reconciliation Schedule" (empty, currently active). Work on the sheet
# that has the data.
$ws = $wb.Worksheets.Item("Statement of Cash Flow")

# --- New column widths for the reconciliation block (S:W) ------------------
$ws.Columns.Item(19).ColumnWidth = 26.83   # S
$ws.Columns.Item(21).ColumnWidth = 2.83    # U
$ws.Columns.Item(22).ColumnWidth = 13.33   # V

# --- Header block ------------------------------------------------------
$ws.Range("U5").Value = "STATEMENT OF RECONCILIATION "
$ws.Range("U5").Font.Bold = $true
$ws.Range("U5").HorizontalAlignment = -4108   # xlCenter

$ws.Range("S7").Value = "Name : Cascade Tiles Co."

$ws.Range("W7").Value = "As Of: Dec 31, 2011"
$ws.Range("W7").HorizontalAlignment = -4152   # xlRight

$ws.Range("V8:W8").Merge()
$ws.Range("V8").Value = "Amount"
$ws.Range("V8:W8").Font.Bold = $true
$ws.Range("V8:W8").HorizontalAlignment = -4108   # xlCenter

$ws.Range("S9").Value = "Particulars"
$ws.Range("S9").Font.Bold = $true
$ws.Range("S9").HorizontalAlignment = -4108

$ws.Range("V9").Value = "Component"
$ws.Range("V9").Font.Bold = $true
$ws.Range("V9").HorizontalAlignment = -4108

$ws.Range("W9").Value = "Total"
$ws.Range("W9").Font.Bold = $true
$ws.Range("W9").HorizontalAlignment = -4108

# --- Body ----------------------------------------------------------------
$ws.Range("S11").Value = "Net Income "
$ws.Range("S11").Font.Bold = $true
$ws.Range("W11").Value = 312
$ws.Range("W11").Font.Bold = $true

$ws.Range("S12").Value = "Depreciation & amortization"
$ws.Range("W12").Value = 45

$ws.Range("S15").Value = "to the net cash flow from operation:"

$ws.Range("S17").Value = "           Increase in Accounts Recievable"
$ws.Range("W17").Value = -110

$ws.Range("S18").Value = "           Increase in Inventory"
$ws.Range("W18").Value = -50

$ws.Range("S19").Value = "           Increase in Accounts Payable"
$ws.Range("W19").Value = 35

$ws.Range("S20").Value = "           Increase in Accrued Wages"
$ws.Range("W20").Value = -10

$ws.Range("S21").Value = "           Increase in Income Taxes Accrued"
$ws.Range("W21").Value = 5

$ws.Range("S14").Value = "Adjustments to reconcile net income"

$ws.Range("S23").Value = "Net cash flow from operation"
$ws.Range("S23").Font.Bold = $true
$ws.Range("W23").Formula = "=SUM(W11:W22)"
$ws.Range("W23").Font.Bold = $true

# --- view settings --------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

# --- drop the empty "Reconciliation Schedule" sheet and rename the
#     remaining sheet to "Reconciliation Statement" ------------------------
$wb.Worksheets.Item("Reconciliation Schedule").Delete()
$ws.Name = "Reconciliation Statement"

$ws.Activate()
$ws.Range("T25").Select()
